$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Airlight Estimator V2 data point (row 6, columns F:I) ---
$ws.Range("F6").Value = "V1.04.1"
$ws.Range("G6").Value = 21.47146
$ws.Range("H6").Value = 0.0072199999999999999
$ws.Range("I6").Value = 0.84968999999999995

# --- Updated Transmission Map Generator rows (13 & 14) ---
$ws.Range("B13").Value = 24.931889999999999
$ws.Range("C13").Value = 0.0032200000000000002
$ws.Range("D13").Value = 0.81723999999999997

$ws.Range("B14").Value = 22.571529999999999
$ws.Range("C14").Value = 0.0055500000000000002
$ws.Range("D14").Value = 0.77834000000000003

# --- New Transmission Map Generator row (15) ---
$ws.Range("A15").Value = "V1.04.3"
$ws.Range("B15").Value = 21.37276
$ws.Range("C15").Value = 0.0073099999999999997
$ws.Range("D15").Value = 0.752

# --- Update the active selection to match the saved view ---
$ws.Range("G9").Select() | Out-Null
